# Daily attendance processing - 2026-01-18 05:16:30
#
# Normalizes the "Recorded By" list (column G) on the active sheet: for any
# cell holding a comma-separated list of recorders, the first and last
# entries are swapped - unless the list is led by the "admin@admin.com"
# account, which always stays first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$recordedByCol = 7   # Column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $v = $cell.Value2

    if ($v -eq $null) {
        continue
    }

    $parts = @($v -split ', ')

    if ($parts.Count -gt 1 -and $parts[0] -ne 'admin@admin.com') {
        $lastIdx = $parts.Count - 1
        $first = $parts[0]
        $last = $parts[$lastIdx]

        $parts[0] = $last
        $parts[$lastIdx] = $first

        $cell.Value = ($parts -join ', ')
    }
}
